$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")

# New rows for integer/bigint regex-fix test data
$ws.Cells.Item(8, 1).Value = "integer neg"
$ws.Cells.Item(8, 2).Value = -2000

$ws.Cells.Item(9, 1).Value = "bigint"
$ws.Cells.Item(9, 2).Value = 100000000000000

$ws.Cells.Item(10, 1).Value = "bigint neg"
$ws.Cells.Item(10, 2).Formula = "=-100000000000000"

# Column B is a bit narrower now that bestFit is no longer forced
$ws.Columns.Item(2).ColumnWidth = 13.5

# general becomes the active sheet/tab, selection on the newly added cell
$ws.Range("B10").Select()
